$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column (2022) to the right of the existing data (column L = 2021).
# New cells mirror the formatting of the neighboring existing cells.

# Header cell M4: year 2022 (reuse style from K4)
$ws.Range("M4").Value = 2022
$ws.Range("K4").Copy()
$ws.Range("M4").PasteSpecial(-4122)

# M5: Small enterprises 2022 value (reuse style from L5)
$ws.Range("M5").Value = 2.2
$ws.Range("L5").Copy()
$ws.Range("M5").PasteSpecial(-4122)

# M6: Medium-sized enterprises 2022 value (reuse style from L6)
$ws.Range("M6").Value = 1.2
$ws.Range("L6").Copy()
$ws.Range("M6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Update selection to match the post-edit active cell
$ws.Range("M10").Select()
